$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ 'B'=0.6194284934348104; 'C'=0.1725593584472165; 'D'=0.0135231891998373; 'E'=0.1282636622104505; 'F'=0.5034795799655143; 'H'=0.07973214163530429; 'I'=0.3408651183409077; 'M'=0.3083193749590549; 'N'=0.9311592783408642; 'O'=1.612172864956676 }
    3 = @{ 'B'=0.5420733954247794; 'C'=0.1519890321154094; 'D'=0.0121326809356006; 'E'=0.1218239920400137; 'F'=0.4948761282576726; 'H'=0.07973214163530429; 'I'=0.3436564154521058; 'M'=0.2744500294856138; 'N'=0.9400192022937546; 'O'=1.596948262432619 }
    4 = @{ 'B'=0.4944770899284379; 'C'=0.1392938856783985; 'D'=0.01127379038398857; 'E'=0.1179836283754412; 'F'=0.4899921034708825; 'H'=0.07973214163530429; 'I'=0.3456378482473887; 'M'=0.2537132601274692; 'N'=0.9458704867437717; 'O'=1.58889822309817 }
    5 = @{ 'B'=0.4750571129636683; 'C'=0.1341043972218756; 'D'=0.01092252354681023; 'E'=0.1164469638191363; 'F'=0.4881018492000919; 'H'=0.07973214163530429; 'I'=0.3465124334721033; 'M'=0.2452777522452934; 'N'=0.9483583729531588; 'O'=1.5859434156595 }
    6 = @{ 'B'=0.4718310159390171; 'C'=0.1332417193912647; 'D'=0.01086412043510165; 'E'=0.1161935045139018; 'F'=0.4877940097165734; 'H'=0.07973214163530429; 'I'=0.3466617082069732; 'M'=0.2438779428478384; 'N'=0.9487777329735039; 'O'=1.585472419308616 }
    7 = @{ 'B'=0.4942152815103213; 'C'=0.1392239633940733; 'D'=0.01126905815489465; 'E'=0.1179627900721201; 'F'=0.4899662060916228; 'H'=0.07973214163530429; 'I'=0.3456493715760871; 'M'=0.2535994355824798; 'N'=0.9459036203849251; 'O'=1.588857055957902 }
    8 = @{ 'B'=0.5927778029784179; 'C'=0.1654802915251992; 'D'=0.01304481692856285; 'E'=0.1260195259336214; 'F'=0.5004302778597491; 'H'=0.07973214163530429; 'I'=0.3417719362567517; 'M'=0.2966289509915541; 'N'=0.9341288392769513; 'O'=1.606653590733885 }
    9 = @{ 'B'=0.7852290036540808; 'C'=0.2164475552764031; 'D'=0.01648562158352718; 'E'=0.1427327368123557; 'F'=0.5241227882272312; 'H'=0.07973214163530429; 'I'=0.3362980872619588; 'M'=0.3814824730978614; 'N'=0.9143006215770342; 'O'=1.651888208096381 }
    10 = @{ 'B'=0.9260818240021536; 'C'=0.2535701757806237; 'D'=0.01898735174653154; 'E'=0.1555876722085117; 'F'=0.5434807257051233; 'H'=0.07973214163530429; 'I'=0.3335844847411025; 'M'=0.4441237877433935; 'N'=0.9017203083722336; 'O'=1.691480876081044 }
    11 = @{ 'B'=0.9900352928331131; 'C'=0.270387185974954; 'D'=0.02011956918151725; 'E'=0.1615648086037282; 'F'=0.5527146132145901; 'H'=0.07973214163530429; 'I'=0.3326361203144046; 'M'=0.4726887606741883; 'N'=0.8964286139695901; 'O'=1.710886183006863 }
    12 = @{ 'B'=1.014234513232168; 'C'=0.2767450911471485; 'D'=0.02054745069010266; 'E'=0.1638470985906935; 'F'=0.556273021384456; 'H'=0.07973214163530429; 'I'=0.3323183020426015; 'M'=0.4835155948606484; 'N'=0.894486790888962; 'O'=1.718435883464963 }
    13 = @{ 'B'=1.00902362224042; 'C'=0.2753762656505501; 'D'=0.02045533762034779; 'E'=0.1633547225914782; 'F'=0.5555039055048923; 'H'=0.07973214163530429; 'I'=0.3323849100495835; 'M'=0.4811834027971003; 'N'=0.8949022380919871; 'O'=1.716800953093241 }
    14 = @{ 'B'=0.992026555977759; 'C'=0.2709104630305035; 'D'=0.02015478872827003; 'E'=0.1617521940440625; 'F'=0.5530061271349638; 'H'=0.07973214163530429; 'I'=0.3326091443882575; 'M'=0.4735792926845761; 'N'=0.896267615734935; 'O'=1.711503261737619 }
    15 = @{ 'B'=0.9816129098123838; 'C'=0.2681736777153105; 'D'=0.01997058047850686; 'E'=0.1607730657366986; 'F'=0.5514842117100898; 'H'=0.07973214163530429; 'I'=0.332751878815273; 'M'=0.4689228435886861; 'N'=0.8971120275688733; 'O'=1.708284516414636 }
    16 = @{ 'B'=0.9218997779911433; 'C'=0.2524697113040304; 'D'=0.01891323909203635; 'E'=0.1551996789620276; 'F'=0.5428858970606569; 'H'=0.07973214163530429; 'I'=0.3336522335804837; 'M'=0.4422583870204164; 'N'=0.9020748111942041; 'O'=1.690240822609553 }
    17 = @{ 'B'=0.8852358523267299; 'C'=0.242817679266949; 'D'=0.01826308211986571; 'E'=0.1518139155644107; 'F'=0.5377208449865236; 'H'=0.07973214163530429; 'I'=0.3342779621945091; 'M'=0.4259182816419411; 'N'=0.9052297649853287; 'O'=1.67952932413553 }
    18 = @{ 'B'=0.8641363936358744; 'C'=0.237259486539358; 'D'=0.01788858162635165; 'E'=0.149878673607617; 'F'=0.5347903038725121; 'H'=0.07973214163530429; 'I'=0.3346647805293301; 'M'=0.4165263939940189; 'N'=0.9070849912036678; 'O'=1.673499552609684 }
    19 = @{ 'B'=0.8569905669904756; 'C'=0.2353764515513319; 'D'=0.01776168896992658; 'E'=0.149225513954903; 'F'=0.5338049813530859; 'H'=0.07973214163530429; 'I'=0.3348003681804173; 'M'=0.4133475813612506; 'N'=0.907720107859852; 'O'=1.67148048418062 }
    20 = @{ 'B'=0.8891399715774924; 'C'=0.2438458392066423; 'D'=0.01833234928201932; 'E'=0.1521730753886601; 'F'=0.538266505296491; 'H'=0.07973214163530429; 'I'=0.3342085653945226; 'M'=0.4276570404212663; 'N'=0.9048897146784398; 'O'=1.680655996881626 }
    21 = @{ 'B'=0.9970195194483722; 'C'=0.2722224597319496; 'D'=0.0202430908969049; 'E'=0.1622223809726719; 'F'=0.5537381080496573; 'H'=0.07973214163530429; 'I'=0.3325421589150146; 'M'=0.4758125358241756; 'N'=0.8958648875948754; 'O'=1.713053851507169 }
    22 = @{ 'B'=1.067416269712794; 'C'=0.2907078952997892; 'D'=0.02148681906189864; 'E'=0.1689003481058577; 'F'=0.5642096319423615; 'H'=0.07973214163530429; 'I'=0.3316939123615619; 'M'=0.5073426153060723; 'N'=0.8903281966029581; 'O'=1.735401616727103 }
    23 = @{ 'B'=1.029854528956946; 'C'=0.2808474657772422; 'D'=0.02082348876317752; 'E'=0.1653260198407409; 'F'=0.5585877795687253; 'H'=0.07973214163530429; 'I'=0.3321245437193383; 'M'=0.4905091542973423; 'N'=0.8932501402443833; 'O'=1.723366512656497 }
    24 = @{ 'B'=0.8873749842159668; 'C'=0.243381036462921; 'D'=0.01830103582821607; 'E'=0.1520106641236012; 'F'=0.5380196910645338; 'H'=0.07973214163530429; 'I'=0.3342398553511998; 'M'=0.4268709405496764; 'N'=0.9050433224582832; 'O'=1.680146228196975 }
    25 = @{ 'B'=0.7332583334279548; 'C'=0.2027159123616684; 'D'=0.0155593298639829; 'E'=0.1381116077358939; 'F'=0.5173719367047696; 'H'=0.07973214163530429; 'I'=0.3375499010877832; 'M'=0.3584755759702958; 'N'=0.9193155903011743; 'O'=1.638538550520082 }
}

foreach ($r in $data.Keys) {
    foreach ($c in $data[$r].Keys) {
        $ws.Range("$c$r").Value = $data[$r][$c]
    }
}
